$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 0.5804443333333333
$ws.Range("N2").Value = 1.741333
$ws.Range("O2").Value = 0.002431273010151717
$ws.Range("P2").Value = 0.002435427107574628
$ws.Range("Q2").Value = 0.08802477011288889
$ws.Range("R2").Value = 0.7922229310160001
$ws.Range("S2").Value = 0.002431273010151717
$ws.Range("T2").Value = 0.002435427107574628

# Row 3
$ws.Range("O3").Value = 0.0004752041289926495
$ws.Range("P3").Value = 0.00047601606752829
$ws.Range("S3").Value = 0.0004752041289926495
$ws.Range("T3").Value = 0.00047601606752829

# Row 4
$ws.Range("M4").Value = 136.1000366666667
$ws.Range("N4").Value = 408.30011
$ws.Range("O4").Value = 0.5700742118164518
$ws.Range("P4").Value = 0.5710482463260632
$ws.Range("Q4").Value = 20.63966129385778
$ws.Range("R4").Value = 185.75695164472
$ws.Range("S4").Value = 0.5700742118164518
$ws.Range("T4").Value = 0.5710482463260632

# Row 5
$ws.Range("M5").Value = 1.221658
$ws.Range("N5").Value = 2.443316
$ws.Range("O5").Value = 0.005117086949542552
$ws.Range("P5").Value = 0.003417220037046797
$ws.Range("Q5").Value = 0.1852652501386667
$ws.Range("R5").Value = 1.111591500832
$ws.Range("S5").Value = 0.005117086949542552
$ws.Range("T5").Value = 0.003417220037046797

# Row 6
$ws.Range("M6").Value = 100.7253213333333
$ws.Range("N6").Value = 302.175964
$ws.Range("O6").Value = 0.4219022240948613
$ws.Range("P6").Value = 0.4226230904617871
$ws.Range("Q6").Value = 15.27506213041423
$ws.Range("R6").Value = 137.475559173728
$ws.Range("S6").Value = 0.4219022240948613
$ws.Range("T6").Value = 0.4226230904617871
